# Auto-generated edit script for hs_meta.xlsx update (3rd commit : project)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "최근 0일 서버의 랭크" -> "최근 1일 서버의 랭크"
$ws.Name = "최근 1일 서버의 랭크"

# Update ranking table cell values
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "56.8%"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "7.0%"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "15,000"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "39.9%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.3%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9,500"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "40.2%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.2%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4,800"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "47.0%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2,100"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "47.3%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1,100"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "53.1%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "780"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "45.9%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.4%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3,100"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "56.1%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.7%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5,900"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "30.5%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.4%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "920"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "48.9%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.2%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "470"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "44.4%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "430"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "26.4%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1,100"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "61.0%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.8%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "12,000"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "58.6%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.2%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "7,000"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "57.7%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.1%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,400"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "35.8%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.2%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,600"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "58.3%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.1%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "6,800"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "57.9%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2,800"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "57.9%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1,800"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "44.8%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.5%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3,200"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.4%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5,200"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "52.4%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.4%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,200"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "61.5%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.7%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3,600"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "47.8%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2,200"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "47.5%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "870"
$ws.Range("B27").Value = "Lightshow Mage"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "44.3%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "410"
$ws.Range("B28").Value = "기계 마법사"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "50.7%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "370"
$ws.Range("B29").Value = "하이랜더 마법사"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "42.8%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "360"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "25.9%"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.9%"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "4,000"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "68.1%"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "13,000"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "66.8%"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.9%"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "8,600"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "61.0%"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.7%"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "3,700"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "53.8%"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7%"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "1,600"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "54.1%"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "390"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "47.0%"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "5,300"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "49.2%"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "5,000"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "57.9%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4,700"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "51.6%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1,300"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "47.6%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "980"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "47.5%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "280"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "32.4%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1,900"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "44.8%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.4%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3,000"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "46.4%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2,200"
$ws.Range("B45").Value = "Mech Rogue"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "57.8%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2,200"
$ws.Range("B46").Value = "Ogre Rogue"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "51.1%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2,100"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "60.7%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "880"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "47.9%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "650"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "48.1%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "490"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "34.6%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.3%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "2,800"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "61.1%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.1%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "4,500"
$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "58.8%"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "0.8%"
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = "1,800"
$ws.Range("C53").NumberFormat = "@"
$ws.Range("C53").Value = "52.6%"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "0.6%"
$ws.Range("E53").NumberFormat = "@"
$ws.Range("E53").Value = "1,300"
$ws.Range("C54").NumberFormat = "@"
$ws.Range("C54").Value = "42.7%"
$ws.Range("E54").NumberFormat = "@"
$ws.Range("E54").Value = "680"
$ws.Range("C55").NumberFormat = "@"
$ws.Range("C55").Value = "36.4%"
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = "0.9%"
$ws.Range("E55").NumberFormat = "@"
$ws.Range("E55").Value = "1,900"
$ws.Range("C56").NumberFormat = "@"
$ws.Range("C56").Value = "57.2%"
$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "3.2%"
$ws.Range("E56").NumberFormat = "@"
$ws.Range("E56").Value = "6,900"
$ws.Range("C57").NumberFormat = "@"
$ws.Range("C57").Value = "52.9%"
$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "1.1%"
$ws.Range("E57").NumberFormat = "@"
$ws.Range("E57").Value = "2,300"
$ws.Range("C58").NumberFormat = "@"
$ws.Range("C58").Value = "48.7%"
$ws.Range("E58").NumberFormat = "@"
$ws.Range("E58").Value = "1,800"
$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "45.2%"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "0.7%"
$ws.Range("E59").NumberFormat = "@"
$ws.Range("E59").Value = "1,600"
$ws.Range("C60").NumberFormat = "@"
$ws.Range("C60").Value = "49.1%"
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = "700"
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "40.6%"
$ws.Range("E61").NumberFormat = "@"
$ws.Range("E61").Value = "3,400"
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").Value = "53.6%"
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "5.7%"
$ws.Range("E62").NumberFormat = "@"
$ws.Range("E62").Value = "12,000"
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "60.6%"
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "1.3%"
$ws.Range("E63").NumberFormat = "@"
$ws.Range("E63").Value = "2,800"
$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "45.3%"
$ws.Range("E64").NumberFormat = "@"
$ws.Range("E64").Value = "1,300"
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "47.2%"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "0.5%"
$ws.Range("E65").NumberFormat = "@"
$ws.Range("E65").Value = "1,000"
$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "48.2%"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "0.4%"
$ws.Range("E66").NumberFormat = "@"
$ws.Range("E66").Value = "810"
$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "37.7%"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = "1.7%"
$ws.Range("E67").NumberFormat = "@"
$ws.Range("E67").Value = "3,700"
